$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds weekly price-report rows, one group of rows per reporting
# date. This commit adds a brand-new weekly group (date 2022-03-25 / serial
# 44578) for "Plátano" and, because the existing groups are *not* sorted by
# date, the new group is inserted right above the group that (in file order)
# used to sit at rows 847-849 -- i.e. the 3 new rows land at 847:849 and
# everything that used to be at row 847 onward is pushed down by 3 rows.

# 1) Insert 3 blank rows at row 850 - this shifts the old rows 850-917 down
#    to 853-920, while rows 847-849 (old data) stay put for now.
$ws.Rows.Item(850).Resize(3).Insert()

# 2) The old content that was in rows 847:849 logically belongs at the new
#    rows 850:852 (it simply moved down by 3 along with everything else).
#    Copy it there now, before we overwrite 847:849 with the new group.
$ws.Range("A847:T849").Copy()
$ws.Range("A850").PasteSpecial()
$excel.CutCopyMode = 0

# 3) Overwrite rows 847:849 with the brand-new weekly group's data.
$ws.Cells.Item(847, 4).Value = 44578
$ws.Cells.Item(847, 13).Value = 650
$ws.Cells.Item(847, 14).Value = 12000
$ws.Cells.Item(847, 15).Value = 12000
$ws.Cells.Item(847, 16).Value = 12000
$ws.Cells.Item(847, 19).Value = 600

$ws.Cells.Item(848, 4).Value = 44578
$ws.Cells.Item(848, 13).Value = 700
$ws.Cells.Item(848, 14).Value = 13000
$ws.Cells.Item(848, 15).Value = 13000
$ws.Cells.Item(848, 16).Value = 13000
$ws.Cells.Item(848, 19).Value = 650

$ws.Cells.Item(849, 4).Value = 44578
$ws.Cells.Item(849, 13).Value = 880
$ws.Cells.Item(849, 14).Value = 15000
$ws.Cells.Item(849, 15).Value = 15000
$ws.Cells.Item(849, 16).Value = 15000
$ws.Cells.Item(849, 19).Value = 750

$ws.Range("A1").Select()
